$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ============================================================
# Section 1: rows 176-188 permutation / odds corrections
# (content among these rows got re-matched to the correct game;
#  row id in column A stays put, everything else moves)
# ============================================================
# Row 176
$ws.Cells.Item(176, 2).Value = 7302795
$ws.Cells.Item(176, 6).Value = "Unin Comercio"
$ws.Cells.Item(176, 7).Value = "Deportivo Garcilaso"
$ws.Cells.Item(176, 9).Value = 2
$ws.Cells.Item(176, 10).Value = "A"
$ws.Cells.Item(176, 11).Value = 2.25
$ws.Cells.Item(176, 12).Value = 3.3
$ws.Cells.Item(176, 13).Value = 2.7
$ws.Cells.Item(176, 14).Value = 1.75
$ws.Cells.Item(176, 15).Value = 3.6
$ws.Cells.Item(176, 16).Value = 4
$ws.Cells.Item(176, 17).Value = -0.5
$ws.Cells.Item(176, 18).Value = 1.8
$ws.Cells.Item(176, 19).Value = 2
$ws.Cells.Item(176, 20).Value = 2.75
$ws.Cells.Item(176, 21).Value = 1.825
$ws.Cells.Item(176, 22).Value = 1.975
$ws.Cells.Item(176, 23).Value = -1
$ws.Cells.Item(176, 25).Value = 3
$ws.Cells.Item(176, 27).Value = 1
$ws.Cells.Item(176, 28).Value = 0.4125
$ws.Cells.Item(176, 29).Value = -0.5

# Row 177
$ws.Cells.Item(177, 2).Value = 7302796
$ws.Cells.Item(177, 6).Value = "Sport Huancayo"
$ws.Cells.Item(177, 7).Value = "Sport Boys"
$ws.Cells.Item(177, 9).Value = 0
$ws.Cells.Item(177, 10).Value = "H"
$ws.Cells.Item(177, 11).Value = 1.727
$ws.Cells.Item(177, 12).Value = 3.75
$ws.Cells.Item(177, 13).Value = 4.333
$ws.Cells.Item(177, 14).Value = 1.25
$ws.Cells.Item(177, 15).Value = 5.25
$ws.Cells.Item(177, 16).Value = 10
$ws.Cells.Item(177, 17).Value = -1.75
$ws.Cells.Item(177, 18).Value = 1.925
$ws.Cells.Item(177, 19).Value = 1.875
$ws.Cells.Item(177, 20).Value = 3
$ws.Cells.Item(177, 21).Value = 1.875
$ws.Cells.Item(177, 22).Value = 1.925
$ws.Cells.Item(177, 23).Value = 0.25
$ws.Cells.Item(177, 25).Value = -1
$ws.Cells.Item(177, 27).Value = 0.875
$ws.Cells.Item(177, 28).Value = -1
$ws.Cells.Item(177, 29).Value = 0.925

# Row 180
$ws.Cells.Item(180, 2).Value = 7384622
$ws.Cells.Item(180, 6).Value = "Deportivo Municipal"
$ws.Cells.Item(180, 7).Value = "Academia Deportiva Cantolao"
$ws.Cells.Item(180, 8).Value = 1
$ws.Cells.Item(180, 9).Value = 2
$ws.Cells.Item(180, 10).Value = "A"
$ws.Cells.Item(180, 11).Value = 1.444
$ws.Cells.Item(180, 12).Value = 4.333
$ws.Cells.Item(180, 13).Value = 7
$ws.Cells.Item(180, 14).Value = 1.5
$ws.Cells.Item(180, 15).Value = 3.75
$ws.Cells.Item(180, 16).Value = 6
$ws.Cells.Item(180, 18).Value = 1.825
$ws.Cells.Item(180, 19).Value = 2.025
$ws.Cells.Item(180, 20).Value = 2.75
$ws.Cells.Item(180, 21).Value = 1.875
$ws.Cells.Item(180, 22).Value = 1.975
$ws.Cells.Item(180, 23).Value = -1
$ws.Cells.Item(180, 25).Value = 5
$ws.Cells.Item(180, 26).Value = -1
$ws.Cells.Item(180, 27).Value = 1.025
$ws.Cells.Item(180, 28).Value = 0.4375
$ws.Cells.Item(180, 29).Value = -0.5

# Row 182
$ws.Cells.Item(182, 2).Value = 7384624
$ws.Cells.Item(182, 6).Value = "Cesar Vallejo"
$ws.Cells.Item(182, 7).Value = "Cusco FC"
$ws.Cells.Item(182, 8).Value = 3
$ws.Cells.Item(182, 9).Value = 1
$ws.Cells.Item(182, 10).Value = "H"
$ws.Cells.Item(182, 11).Value = 2
$ws.Cells.Item(182, 12).Value = 3.4
$ws.Cells.Item(182, 13).Value = 3.5
$ws.Cells.Item(182, 14).Value = 1.45
$ws.Cells.Item(182, 15).Value = 4.2
$ws.Cells.Item(182, 16).Value = 6.5
$ws.Cells.Item(182, 18).Value = 1.75
$ws.Cells.Item(182, 19).Value = 2.05
$ws.Cells.Item(182, 20).Value = 2.5
$ws.Cells.Item(182, 21).Value = 1.95
$ws.Cells.Item(182, 22).Value = 1.85
$ws.Cells.Item(182, 23).Value = 0.45
$ws.Cells.Item(182, 25).Value = -1
$ws.Cells.Item(182, 26).Value = 0.75
$ws.Cells.Item(182, 27).Value = -1
$ws.Cells.Item(182, 28).Value = 0.95
$ws.Cells.Item(182, 29).Value = -1

# Row 183
$ws.Cells.Item(183, 2).Value = 7384630
$ws.Cells.Item(183, 6).Value = "Atletico Grau"
$ws.Cells.Item(183, 7).Value = "Unin Comercio"
$ws.Cells.Item(183, 9).Value = 1
$ws.Cells.Item(183, 10).Value = "A"
$ws.Cells.Item(183, 11).Value = 2.8
$ws.Cells.Item(183, 12).Value = 3.4
$ws.Cells.Item(183, 13).Value = 2.15
$ws.Cells.Item(183, 14).Value = 1.75
$ws.Cells.Item(183, 15).Value = 3.6
$ws.Cells.Item(183, 16).Value = 3.8
$ws.Cells.Item(183, 17).Value = -0.75
$ws.Cells.Item(183, 18).Value = 2
$ws.Cells.Item(183, 19).Value = 1.8
$ws.Cells.Item(183, 20).Value = 3
$ws.Cells.Item(183, 21).Value = 1.85
$ws.Cells.Item(183, 22).Value = 1.95
$ws.Cells.Item(183, 24).Value = -1
$ws.Cells.Item(183, 25).Value = 2.8
$ws.Cells.Item(183, 27).Value = 0.8
$ws.Cells.Item(183, 29).Value = 0.95

# Row 184
$ws.Cells.Item(184, 2).Value = 7384625
$ws.Cells.Item(184, 6).Value = "AD Tarma"
$ws.Cells.Item(184, 7).Value = "Carlos Manucci"
$ws.Cells.Item(184, 8).Value = 0
$ws.Cells.Item(184, 10).Value = "D"
$ws.Cells.Item(184, 11).Value = 1.5
$ws.Cells.Item(184, 12).Value = 3.75
$ws.Cells.Item(184, 13).Value = 7
$ws.Cells.Item(184, 14).Value = 1.363
$ws.Cells.Item(184, 15).Value = 4.333
$ws.Cells.Item(184, 16).Value = 9.5
$ws.Cells.Item(184, 17).Value = -1.25
$ws.Cells.Item(184, 18).Value = 1.875
$ws.Cells.Item(184, 19).Value = 1.925
$ws.Cells.Item(184, 20).Value = 2.5
$ws.Cells.Item(184, 21).Value = 1.8
$ws.Cells.Item(184, 22).Value = 2
$ws.Cells.Item(184, 23).Value = -1
$ws.Cells.Item(184, 24).Value = 3.333
$ws.Cells.Item(184, 26).Value = -1
$ws.Cells.Item(184, 27).Value = 0.925
$ws.Cells.Item(184, 28).Value = -1
$ws.Cells.Item(184, 29).Value = 1

# Row 185
$ws.Cells.Item(185, 2).Value = 7384628
$ws.Cells.Item(185, 6).Value = "Deportivo Binacional"
$ws.Cells.Item(185, 7).Value = "FBC Melgar"
$ws.Cells.Item(185, 8).Value = 1
$ws.Cells.Item(185, 9).Value = 2
$ws.Cells.Item(185, 11).Value = 2.75
$ws.Cells.Item(185, 13).Value = 2.375
$ws.Cells.Item(185, 14).Value = 3.3
$ws.Cells.Item(185, 15).Value = 3.6
$ws.Cells.Item(185, 16).Value = 2
$ws.Cells.Item(185, 17).Value = 0.5
$ws.Cells.Item(185, 18).Value = 1.8
$ws.Cells.Item(185, 19).Value = 2
$ws.Cells.Item(185, 20).Value = 2.75
$ws.Cells.Item(185, 21).Value = 1.975
$ws.Cells.Item(185, 22).Value = 1.875
$ws.Cells.Item(185, 25).Value = 1
$ws.Cells.Item(185, 27).Value = 1
$ws.Cells.Item(185, 28).Value = 0.4875
$ws.Cells.Item(185, 29).Value = -0.5

# Row 186
$ws.Cells.Item(186, 2).Value = 7384627
$ws.Cells.Item(186, 6).Value = "Universitario de Deportes"
$ws.Cells.Item(186, 7).Value = "Sport Huancayo"
$ws.Cells.Item(186, 8).Value = 2
$ws.Cells.Item(186, 9).Value = 0
$ws.Cells.Item(186, 10).Value = "H"
$ws.Cells.Item(186, 11).Value = 1.25
$ws.Cells.Item(186, 12).Value = 5
$ws.Cells.Item(186, 13).Value = 12
$ws.Cells.Item(186, 14).Value = 1.181
$ws.Cells.Item(186, 15).Value = 6
$ws.Cells.Item(186, 16).Value = 13
$ws.Cells.Item(186, 17).Value = -1.75
$ws.Cells.Item(186, 21).Value = 1.85
$ws.Cells.Item(186, 22).Value = 1.95
$ws.Cells.Item(186, 23).Value = 0.181
$ws.Cells.Item(186, 25).Value = -1
$ws.Cells.Item(186, 26).Value = 0.4
$ws.Cells.Item(186, 27).Value = -0.5
$ws.Cells.Item(186, 28).Value = -1
$ws.Cells.Item(186, 29).Value = 0.95

# Row 187
$ws.Cells.Item(187, 2).Value = 7384626
$ws.Cells.Item(187, 6).Value = "Sporting Cristal"
$ws.Cells.Item(187, 7).Value = "Alianza Atletico"
$ws.Cells.Item(187, 8).Value = 3
$ws.Cells.Item(187, 9).Value = 0
$ws.Cells.Item(187, 10).Value = "H"
$ws.Cells.Item(187, 11).Value = 1.3
$ws.Cells.Item(187, 12).Value = 5
$ws.Cells.Item(187, 13).Value = 9
$ws.Cells.Item(187, 14).Value = 1.166
$ws.Cells.Item(187, 15).Value = 6.5
$ws.Cells.Item(187, 16).Value = 13
$ws.Cells.Item(187, 17).Value = -2
$ws.Cells.Item(187, 18).Value = 1.85
$ws.Cells.Item(187, 19).Value = 1.95
$ws.Cells.Item(187, 20).Value = 3.25
$ws.Cells.Item(187, 21).Value = 2
$ws.Cells.Item(187, 22).Value = 1.8
$ws.Cells.Item(187, 23).Value = 0.1659999999999999
$ws.Cells.Item(187, 25).Value = -1
$ws.Cells.Item(187, 26).Value = 0.8500000000000001
$ws.Cells.Item(187, 27).Value = -1
$ws.Cells.Item(187, 28).Value = -0.5
$ws.Cells.Item(187, 29).Value = 0.4

# Row 188
$ws.Cells.Item(188, 2).Value = 7384629
$ws.Cells.Item(188, 6).Value = "Deportivo Garcilaso"
$ws.Cells.Item(188, 7).Value = "Alianza Lima"
$ws.Cells.Item(188, 8).Value = 0
$ws.Cells.Item(188, 9).Value = 1
$ws.Cells.Item(188, 10).Value = "A"
$ws.Cells.Item(188, 11).Value = 2.625
$ws.Cells.Item(188, 12).Value = 3.3
$ws.Cells.Item(188, 13).Value = 2.5
$ws.Cells.Item(188, 14).Value = 2.7
$ws.Cells.Item(188, 15).Value = 3.4
$ws.Cells.Item(188, 16).Value = 2.375
$ws.Cells.Item(188, 17).Value = 0
$ws.Cells.Item(188, 18).Value = 2.025
$ws.Cells.Item(188, 19).Value = 1.775
$ws.Cells.Item(188, 20).Value = 2.25
$ws.Cells.Item(188, 21).Value = 1.825
$ws.Cells.Item(188, 22).Value = 1.975
$ws.Cells.Item(188, 23).Value = -1
$ws.Cells.Item(188, 25).Value = 1.375
$ws.Cells.Item(188, 26).Value = -1
$ws.Cells.Item(188, 27).Value = 0.7749999999999999
$ws.Cells.Item(188, 29).Value = 0.9750000000000001

# ============================================================
# Section 2: rows 281-289 (new results + new fixtures)
# ============================================================

# Make room: push the old last row (8071407, currently row 283)
# down to row 284, and stamp out template rows 285-289 with correct
# styling (column A / E formats) by copying that same template row.
$ws.Range("A283:AC283").Copy($ws.Range("A284:AC284"))
$ws.Range("A283:AC283").Copy($ws.Range("A285:AC285"))
$ws.Range("A283:AC283").Copy($ws.Range("A286:AC286"))
$ws.Range("A283:AC283").Copy($ws.Range("A287:AC287"))
$ws.Range("A283:AC283").Copy($ws.Range("A288:AC288"))
$ws.Range("A283:AC283").Copy($ws.Range("A289:AC289"))

# Clear the phantom H/I/J/AB/AC cells the copy brought along on 284-289
$ws.Range("H284:J289").ClearContents()
$ws.Range("AB284:AC289").ClearContents()

# Row 281
$ws.Cells.Item(281, 1).Value = 279
$ws.Cells.Item(281, 2).Value = 8042070
$ws.Cells.Item(281, 3).Value = "Peru Liga 1"
$ws.Cells.Item(281, 4).Value = "Peru Liga 1"
$ws.Cells.Item(281, 5).Value = 45394.70833333334
$ws.Cells.Item(281, 6).Value = "Sport Huancayo"
$ws.Cells.Item(281, 7).Value = "Cienciano"
$ws.Cells.Item(281, 8).Value = 1
$ws.Cells.Item(281, 9).Value = 2
$ws.Cells.Item(281, 10).Value = "A"
$ws.Cells.Item(281, 11).Value = 1.8
$ws.Cells.Item(281, 12).Value = 3.5
$ws.Cells.Item(281, 13).Value = 4.333
$ws.Cells.Item(281, 14).Value = 2.2
$ws.Cells.Item(281, 15).Value = 3.3
$ws.Cells.Item(281, 16).Value = 3.2
$ws.Cells.Item(281, 17).Value = -0.25
$ws.Cells.Item(281, 18).Value = 1.925
$ws.Cells.Item(281, 19).Value = 1.875
$ws.Cells.Item(281, 20).Value = 2.5
$ws.Cells.Item(281, 21).Value = 1.875
$ws.Cells.Item(281, 22).Value = 1.925
$ws.Cells.Item(281, 23).Value = -1
$ws.Cells.Item(281, 24).Value = -1
$ws.Cells.Item(281, 25).Value = 2.2
$ws.Cells.Item(281, 26).Value = -1
$ws.Cells.Item(281, 27).Value = 0.875
$ws.Cells.Item(281, 28).Value = 0.875
$ws.Cells.Item(281, 29).Value = -1

# Row 282
$ws.Cells.Item(282, 1).Value = 280
$ws.Cells.Item(282, 2).Value = 8042071
$ws.Cells.Item(282, 3).Value = "Peru Liga 1"
$ws.Cells.Item(282, 4).Value = "Peru Liga 1"
$ws.Cells.Item(282, 5).Value = 45395.6875
$ws.Cells.Item(282, 6).Value = "Alianza Atletico"
$ws.Cells.Item(282, 7).Value = "Union Comercio"
$ws.Cells.Item(282, 11).Value = 1.727
$ws.Cells.Item(282, 12).Value = 3.5
$ws.Cells.Item(282, 13).Value = 5
$ws.Cells.Item(282, 14).Value = 1.615
$ws.Cells.Item(282, 15).Value = 3.6
$ws.Cells.Item(282, 16).Value = 6
$ws.Cells.Item(282, 17).Value = -0.75
$ws.Cells.Item(282, 18).Value = 1.825
$ws.Cells.Item(282, 19).Value = 2.025
$ws.Cells.Item(282, 20).Value = 2.25
$ws.Cells.Item(282, 21).Value = 1.925
$ws.Cells.Item(282, 22).Value = 1.925
$ws.Cells.Item(282, 23).Value = 0
$ws.Cells.Item(282, 24).Value = 0
$ws.Cells.Item(282, 25).Value = 0
$ws.Cells.Item(282, 26).Value = 0
$ws.Cells.Item(282, 27).Value = 0

# Row 283
$ws.Cells.Item(283, 1).Value = 281
$ws.Cells.Item(283, 2).Value = 8042072
$ws.Cells.Item(283, 3).Value = "Peru Liga 1"
$ws.Cells.Item(283, 4).Value = "Peru Liga 1"
$ws.Cells.Item(283, 5).Value = 45395.79166666666
$ws.Cells.Item(283, 6).Value = "FBC Melgar"
$ws.Cells.Item(283, 7).Value = "Carlos Manucci"
$ws.Cells.Item(283, 11).Value = 1.25
$ws.Cells.Item(283, 12).Value = 5.5
$ws.Cells.Item(283, 13).Value = 12
$ws.Cells.Item(283, 14).Value = 1.222
$ws.Cells.Item(283, 15).Value = 6
$ws.Cells.Item(283, 16).Value = 13
$ws.Cells.Item(283, 17).Value = -1.75
$ws.Cells.Item(283, 18).Value = 1.875
$ws.Cells.Item(283, 19).Value = 1.975
$ws.Cells.Item(283, 20).Value = 3
$ws.Cells.Item(283, 21).Value = 1.85
$ws.Cells.Item(283, 22).Value = 2
$ws.Cells.Item(283, 23).Value = 0
$ws.Cells.Item(283, 24).Value = 0
$ws.Cells.Item(283, 25).Value = 0
$ws.Cells.Item(283, 26).Value = 0
$ws.Cells.Item(283, 27).Value = 0

# Row 284
$ws.Cells.Item(284, 1).Value = 282
$ws.Cells.Item(284, 2).Value = 8071407
$ws.Cells.Item(284, 3).Value = "Peru Liga 1"
$ws.Cells.Item(284, 4).Value = "Peru Liga 1"
$ws.Cells.Item(284, 5).Value = 45395.91666666666
$ws.Cells.Item(284, 6).Value = "Sport Boys"
$ws.Cells.Item(284, 7).Value = "Universitario de Deportes"
$ws.Cells.Item(284, 11).Value = 6
$ws.Cells.Item(284, 12).Value = 3.75
$ws.Cells.Item(284, 13).Value = 1.571
$ws.Cells.Item(284, 14).Value = 5.75
$ws.Cells.Item(284, 15).Value = 3.6
$ws.Cells.Item(284, 16).Value = 1.615
$ws.Cells.Item(284, 17).Value = 0.75
$ws.Cells.Item(284, 18).Value = 2.05
$ws.Cells.Item(284, 19).Value = 1.8
$ws.Cells.Item(284, 20).Value = 2.25
$ws.Cells.Item(284, 21).Value = 1.825
$ws.Cells.Item(284, 22).Value = 2.025
$ws.Cells.Item(284, 23).Value = 0
$ws.Cells.Item(284, 24).Value = 0
$ws.Cells.Item(284, 25).Value = 0
$ws.Cells.Item(284, 26).Value = 0
$ws.Cells.Item(284, 27).Value = 0

# Row 285
$ws.Cells.Item(285, 1).Value = 283
$ws.Cells.Item(285, 2).Value = 8042215
$ws.Cells.Item(285, 3).Value = "Peru Liga 1"
$ws.Cells.Item(285, 4).Value = "Peru Liga 1"
$ws.Cells.Item(285, 5).Value = 45396.625
$ws.Cells.Item(285, 6).Value = "Cesar Vallejo"
$ws.Cells.Item(285, 7).Value = "Comerciantes Unidos"
$ws.Cells.Item(285, 11).Value = 1.8
$ws.Cells.Item(285, 12).Value = 3.5
$ws.Cells.Item(285, 13).Value = 4.333
$ws.Cells.Item(285, 14).Value = 1.45
$ws.Cells.Item(285, 15).Value = 4
$ws.Cells.Item(285, 16).Value = 7.5
$ws.Cells.Item(285, 17).Value = -1.25
$ws.Cells.Item(285, 18).Value = 2.025
$ws.Cells.Item(285, 19).Value = 1.825
$ws.Cells.Item(285, 20).Value = 2.75
$ws.Cells.Item(285, 21).Value = 2
$ws.Cells.Item(285, 22).Value = 1.85
$ws.Cells.Item(285, 23).Value = 0
$ws.Cells.Item(285, 24).Value = 0
$ws.Cells.Item(285, 25).Value = 0
$ws.Cells.Item(285, 26).Value = 0
$ws.Cells.Item(285, 27).Value = 0

# Row 286
$ws.Cells.Item(286, 1).Value = 284
$ws.Cells.Item(286, 2).Value = 8042073
$ws.Cells.Item(286, 3).Value = "Peru Liga 1"
$ws.Cells.Item(286, 4).Value = "Peru Liga 1"
$ws.Cells.Item(286, 5).Value = 45396.79166666666
$ws.Cells.Item(286, 6).Value = "Cusco FC"
$ws.Cells.Item(286, 7).Value = "AD Tarma"
$ws.Cells.Item(286, 11).Value = 1.909
$ws.Cells.Item(286, 12).Value = 3.6
$ws.Cells.Item(286, 13).Value = 3.6
$ws.Cells.Item(286, 14).Value = 2
$ws.Cells.Item(286, 15).Value = 3.6
$ws.Cells.Item(286, 16).Value = 3.3
$ws.Cells.Item(286, 17).Value = -0.5
$ws.Cells.Item(286, 18).Value = 2.05
$ws.Cells.Item(286, 19).Value = 1.8
$ws.Cells.Item(286, 20).Value = 2.5
$ws.Cells.Item(286, 21).Value = 1.925
$ws.Cells.Item(286, 22).Value = 1.925
$ws.Cells.Item(286, 23).Value = 0
$ws.Cells.Item(286, 24).Value = 0
$ws.Cells.Item(286, 25).Value = 0
$ws.Cells.Item(286, 26).Value = 0
$ws.Cells.Item(286, 27).Value = 0

# Row 287
$ws.Cells.Item(287, 1).Value = 285
$ws.Cells.Item(287, 2).Value = 8042219
$ws.Cells.Item(287, 3).Value = "Peru Liga 1"
$ws.Cells.Item(287, 4).Value = "Peru Liga 1"
$ws.Cells.Item(287, 5).Value = 45396.89583333334
$ws.Cells.Item(287, 6).Value = "Alianza Lima"
$ws.Cells.Item(287, 7).Value = "Atletico Grau"
$ws.Cells.Item(287, 11).Value = 1.363
$ws.Cells.Item(287, 12).Value = 4.5
$ws.Cells.Item(287, 13).Value = 9
$ws.Cells.Item(287, 14).Value = 1.4
$ws.Cells.Item(287, 15).Value = 4.333
$ws.Cells.Item(287, 16).Value = 8.5
$ws.Cells.Item(287, 17).Value = -1.25
$ws.Cells.Item(287, 18).Value = 1.875
$ws.Cells.Item(287, 19).Value = 1.975
$ws.Cells.Item(287, 20).Value = 2.75
$ws.Cells.Item(287, 21).Value = 2
$ws.Cells.Item(287, 22).Value = 1.85
$ws.Cells.Item(287, 23).Value = 0
$ws.Cells.Item(287, 24).Value = 0
$ws.Cells.Item(287, 25).Value = 0
$ws.Cells.Item(287, 26).Value = 0
$ws.Cells.Item(287, 27).Value = 0

# Row 288
$ws.Cells.Item(288, 1).Value = 286
$ws.Cells.Item(288, 2).Value = 8042075
$ws.Cells.Item(288, 3).Value = "Peru Liga 1"
$ws.Cells.Item(288, 4).Value = "Peru Liga 1"
$ws.Cells.Item(288, 5).Value = 45397.70833333334
$ws.Cells.Item(288, 6).Value = "UTC Cajamarca"
$ws.Cells.Item(288, 7).Value = "CD Los Chankas"
$ws.Cells.Item(288, 11).Value = 1.909
$ws.Cells.Item(288, 12).Value = 3.5
$ws.Cells.Item(288, 13).Value = 3.75
$ws.Cells.Item(288, 14).Value = 1.909
$ws.Cells.Item(288, 15).Value = 3.5
$ws.Cells.Item(288, 16).Value = 3.8
$ws.Cells.Item(288, 17).Value = -0.5
$ws.Cells.Item(288, 18).Value = 1.95
$ws.Cells.Item(288, 19).Value = 1.9
$ws.Cells.Item(288, 20).Value = 2.5
$ws.Cells.Item(288, 21).Value = 1.925
$ws.Cells.Item(288, 22).Value = 1.925
$ws.Cells.Item(288, 23).Value = 0
$ws.Cells.Item(288, 24).Value = 0
$ws.Cells.Item(288, 25).Value = 0
$ws.Cells.Item(288, 26).Value = 0
$ws.Cells.Item(288, 27).Value = 0

# Row 289
$ws.Cells.Item(289, 1).Value = 287
$ws.Cells.Item(289, 2).Value = 8071422
$ws.Cells.Item(289, 3).Value = "Peru Liga 1"
$ws.Cells.Item(289, 4).Value = "Peru Liga 1"
$ws.Cells.Item(289, 5).Value = 45397.89583333334
$ws.Cells.Item(289, 6).Value = "Deportivo Garcilaso"
$ws.Cells.Item(289, 7).Value = "Sporting Cristal"
$ws.Cells.Item(289, 11).Value = 4
$ws.Cells.Item(289, 12).Value = 3.6
$ws.Cells.Item(289, 13).Value = 1.833
$ws.Cells.Item(289, 14).Value = 3.75
$ws.Cells.Item(289, 15).Value = 3.6
$ws.Cells.Item(289, 16).Value = 1.909
$ws.Cells.Item(289, 17).Value = 0.5
$ws.Cells.Item(289, 18).Value = 1.9
$ws.Cells.Item(289, 19).Value = 1.95
$ws.Cells.Item(289, 20).Value = 2.5
$ws.Cells.Item(289, 21).Value = 1.85
$ws.Cells.Item(289, 22).Value = 2
$ws.Cells.Item(289, 23).Value = 0
$ws.Cells.Item(289, 24).Value = 0
$ws.Cells.Item(289, 25).Value = 0
$ws.Cells.Item(289, 26).Value = 0
$ws.Cells.Item(289, 27).Value = 0
